$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated symbol list values (price in column D, 1h volume change % in column E).
# A leading apostrophe forces Excel's text entry (quote-prefix) so cells that look
# like numbers/percentages are stored as literal text, matching the original inlineStr cells.

$ws.Range("D2").Value = "'296.94"
$ws.Range("E2").Value = "'-2.11%"

$ws.Range("D3").Value = "'31.44"
$ws.Range("E3").Value = "'-1.58%"

$ws.Range("D4").Value = "'5.109"
$ws.Range("E4").Value = "'-2.17%"

$ws.Range("D5").Value = "'0.07332"
$ws.Range("E5").Value = "'-1.42%"

$ws.Range("D6").Value = "'7.723"
$ws.Range("E6").Value = "'-1.33%"

$ws.Range("D7").Value = "'1.731"
$ws.Range("E7").Value = "'18.26%"

$ws.Range("E8").Value = "'-0.27%"

$ws.Range("D9").Value = "'0.9259"
$ws.Range("E9").Value = "'2.26%"

$ws.Range("D10").Value = "'0.1675"
$ws.Range("E10").Value = "'-0.29%"

$ws.Range("D11").Value = "'0.07010"
$ws.Range("E11").Value = "'-5.76%"

$ws.Range("D12").Value = "'0.07971"
$ws.Range("E12").Value = "'-0.68%"

$ws.Range("D13").Value = "'0.02995"
$ws.Range("E13").Value = "'-1.34%"

$ws.Range("D14").Value = "'0.09901"
$ws.Range("E14").Value = "'-0.09%"

$ws.Range("D15").Value = "'0.001500"
$ws.Range("E15").Value = "'0.12%"

$ws.Range("D16").Value = "'0.006253"
$ws.Range("E16").Value = "'-2.36%"

$ws.Range("D17").Value = "'3.451"
$ws.Range("E17").Value = "'-1.09%"

$ws.Range("D18").Value = "'2.221"
$ws.Range("E18").Value = "'-0.43%"

$ws.Range("D19").Value = "'0.3270"
$ws.Range("E19").Value = "'-2.03%"

$ws.Range("D20").Value = "'0.1332"
$ws.Range("E20").Value = "'-0.08%"

$ws.Range("D21").Value = "'4.554"
$ws.Range("E21").Value = "'1.41%"

$ws.Range("D22").Value = "'0.04643"
$ws.Range("E22").Value = "'2.27%"

$ws.Range("D23").Value = "'0.1582"
$ws.Range("E23").Value = "'-4.04%"

$ws.Range("D24").Value = "'0.001221"
$ws.Range("E24").Value = "'0.51%"

$ws.Range("D25").Value = "'0.004746"
$ws.Range("E25").Value = "'7.06%"

$ws.Range("D26").Value = "'0.0001298"

$ws.Range("D27").Value = "'0.0001873"
$ws.Range("E27").Value = "'7.75%"

$ws.Range("D39").Value = "'0.01702"
$ws.Range("E39").Value = "'2.22%"

$ws.Range("D40").Value = "'0.04436"
$ws.Range("E40").Value = "'-1.37%"

$ws.Range("D41").Value = "'0.007212"

$ws.Range("E42").Value = "'-1.33%"

$ws.Range("D43").Value = "'0.002217"
$ws.Range("E43").Value = "'-7.10%"

$ws.Range("D44").Value = "'0.01095"
$ws.Range("E44").Value = "'-21.77%"

$ws.Range("D45").Value = "'0.00006016"
$ws.Range("E45").Value = "'-2.55%"

$ws.Range("E46").Value = "'-21.32%"

$ws.Range("D47").Value = "'0.7211"
$ws.Range("E47").Value = "'2.03%"
